$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the row that held "M3-10mm nylon Screw" (row 42). This shifts all
# subsequent rows up by one, matching the new row count (dimension O68->O67,
# table B4:F60->B4:F59).
$ws.Rows.Item(42).Delete()

# After the shift, the old "M3 nylon nut" row now occupies row 42. Its
# Components/Units values are cleared, leaving the row blank (D/E/F already
# carry on as before: D blank, E "any store", F blank).
$ws.Range("B42").Value = ""
$ws.Range("C42").Value = ""

# Quantity tweaks made alongside the row removal.
$ws.Range("C32").Value = 42
$ws.Range("C47").Value = 2
$ws.Range("C51").Value = 4
$ws.Range("C53").Value = 1

# Restore the view's selection state saved with the workbook.
$ws.Range("C40").Select()
